# Apply the "department" zone additions + active-tab switch (combaine <- tractor)
$wb = $excel.ActiveWorkbook

# --- 1. Add the six new zone/abbreviation rows to the "department" sheet ---
$dept = $wb.Worksheets.Item("department")

$newRows = @(
    @("Роїще", "Ро"),
    @("Новий Білоус", "НБ"),
    @("Рудка", "Ру"),
    @("Пльохів", "Пл"),
    @("Жукотки", "Жу"),
    @("Велика Вісь", "ВВ")
)

$startRow = 8
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $dept.Cells.Item($r, 1).Value = $newRows[$i][0]
    $dept.Cells.Item($r, 2).Value = $newRows[$i][1]
}

$dept.Activate()
$dept.Range("L12").Select()

# --- 2. Switch the selected tab from "tractor" to "combaine" ---
$tractor = $wb.Worksheets.Item("tractor")
$tractor.Activate()
$tractor.Range("G25").Select()

$combaine = $wb.Worksheets.Item("combaine")
$combaine.Activate()
$combaine.Range("C28").Select()

Write-Output "done"
